# Add new daily COVID death-cumulative records (rows 144-159) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Date (serial), DeathCovid, DeathWithCovid, Total
$newData = @(
    @(44263, 8037, 1601,  9638),
    @(44264, 8146, 1614,  9760),
    @(44265, 8244, 1632,  9876),
    @(44266, 8346, 1639,  9985),
    @(44267, 8440, 1649, 10089),
    @(44268, 8528, 1657, 10185),
    @(44269, 8605, 1666, 10271),
    @(44270, 8669, 1673, 10342),
    @(44271, 8738, 1687, 10425),
    @(44272, 8814, 1694, 10508),
    @(44273, 8894, 1707, 10601),
    @(44274, 8978, 1710, 10688),
    @(44275, 9044, 1710, 10754),
    @(44276, 9104, 1716, 10820),
    @(44277, 9190, 1718, 10908),
    @(44278, 9260, 1718, 10978)
)

$row = 144
foreach ($rec in $newData) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row = $row + 1
}

$lastRow = $row - 1

# Update the sheet view to mirror the author's scroll/selection position.
$ws.Activate()
[void]$ws.Range("A" + $lastRow).Select()
$excel.ActiveWindow.ScrollRow = 113
$excel.ActiveWindow.ScrollColumn = 1
